$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the header style/format from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I and J data columns for rows 2-17
$values = @{
    2  = @(9, 9)
    3  = @(5, 6)
    4  = @(9, 9)
    5  = @(3, 4)
    6  = @(7, 7)
    7  = @(8, 8)
    8  = @(6, 6)
    9  = @(6, 6)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(5, 7)
    14 = @(1, 3)
    15 = @(4, 6)
    16 = @(1, 2)
    17 = @(3, 4)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Range("I$r").Value = $pair[0]
    $ws.Range("J$r").Value = $pair[1]
}
